# Harmonize similar tags to be the same.
#
# The "SwateTemplateMetadata" sheet holds a #TAGS list block (rows 12-14):
#   Row 12 : Tags                         | Plant | growth chamber | growth protocol
#   Row 13 : Tags Term Accession Number   |       NCIT_C14258 URL  |   EFO_0003789 URL
#   Row 14 : Tags Term Source REF         |       NCIT             |   EFO
#
# "growth protocol" is renamed to "growth" so it lines up with the already
# existing GO term for "growth" (GO:0040007) instead of the EFO protocol term,
# and the accession numbers are harmonized to use compact CURIE notation
# (e.g. "NCIT:C14258") instead of full purl.obolibrary.org URLs. Because the
# CURIE already carries the ontology prefix, the separate "Term Source REF"
# values for these two tags are no longer needed and are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SwateTemplateMetadata")

# Row 12: rename the third tag from "growth protocol" to "growth"
$ws.Cells.Item(12, 4).Value = "growth"

# Row 13: Tags Term Accession Number -> compact CURIEs
$ws.Cells.Item(13, 2).Value = "NCIT:C14258"
$ws.Cells.Item(13, 4).Value = "GO:0040007"

# Row 14: Tags Term Source REF -> no longer needed, clear them
$ws.Cells.Item(14, 2).ClearContents()
$ws.Cells.Item(14, 4).ClearContents()
